$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Trt"
$ws.Range("B1").Value = "Trt"
$ws.Range("C1").Value = "Dry wt"

$ws.Range("B6").Select()
